$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Stamp the formatting for the 12 new rows (284-295) by copying it from
#    existing rows that already carry the right alternating style:
#      - even data rows use the style quadruple seen on row 2  (A=4, B..M=5, N=6)
#      - odd  data rows use the style quadruple seen on row 3  (A=7, B..M=8, N=9)
#    Row 295 is the new LAST row of the table, so it must carry the special
#    "closing" style that row 283 currently has (A=15, B..L=16, M/N=19).
# ---------------------------------------------------------------------------

$ws.Range("A2:N2").Copy()
$ws.Range("A284:N284").PasteSpecial(-4122)
$ws.Range("A286:N286").PasteSpecial(-4122)
$ws.Range("A288:N288").PasteSpecial(-4122)
$ws.Range("A290:N290").PasteSpecial(-4122)
$ws.Range("A292:N292").PasteSpecial(-4122)
$ws.Range("A294:N294").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A3:N3").Copy()
$ws.Range("A285:N285").PasteSpecial(-4122)
$ws.Range("A287:N287").PasteSpecial(-4122)
$ws.Range("A289:N289").PasteSpecial(-4122)
$ws.Range("A291:N291").PasteSpecial(-4122)
$ws.Range("A293:N293").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Capture the "last row" closing style from row 283 before we repurpose it,
# then stamp it onto the new last row (295).
$ws.Range("A283:N283").Copy()
$ws.Range("A295:N295").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Fill in the survey-response data for rows 284-295.
# ---------------------------------------------------------------------------

$ws.Range("A284").Value = 45599.5927162963
$ws.Range("B284").Value = "bevery2685@gmail.com"
$ws.Range("C284").Value = "반도체디스플레이스쿨"
$ws.Range("D284").Value = 20243354
$ws.Range("E284").Value = "조영태"
$ws.Range("F284").Value = "대한민국"
$ws.Range("G284").Value = "취업자 / 15세 이상 인구"
$ws.Range("H284").Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Range("I284").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J284").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K284").Value = "`"19.1%`""
$ws.Range("L284").Value = "Black"
$ws.Range("N284").Value = "나는 사후 장기기증에 참여할 뜻이 있다"

$ws.Range("A285").Value = 45599.609972349535
$ws.Range("B285").Value = "qkfdmltls@naver.com"
$ws.Range("C285").Value = "사회복지과"
$ws.Range("D285").Value = 20227022
$ws.Range("E285").Value = "조혜람"
$ws.Range("F285").Value = "대한민국"
$ws.Range("G285").Value = "취업자 / 경제활동인구"
$ws.Range("H285").Value = "조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자"
$ws.Range("I285").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J285").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K285").Value = "`"15%`""
$ws.Range("L285").Value = "Red"
$ws.Range("M285").Value = "나는 사후 장기기증에 참여할 뜻이 없다"

$ws.Range("A286").Value = 45599.613347199076
$ws.Range("B286").Value = "tngusvhs@gmail.com"
$ws.Range("C286").Value = "생명과학과"
$ws.Range("D286").Value = 20243529
$ws.Range("E286").Value = "이수현"
$ws.Range("F286").Value = "대한민국"
$ws.Range("G286").Value = "취업자 / 15세 이상 인구"
$ws.Range("H286").Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Range("I286").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J286").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K286").Value = "`"10%`""
$ws.Range("L286").Value = "Black"

$ws.Range("A287").Value = 45599.61532965278
$ws.Range("B287").Value = "detect5641@gmail.com"
$ws.Range("C287").Value = "소프트웨어학부"
$ws.Range("D287").Value = 20245230
$ws.Range("E287").Value = "이주성"
$ws.Range("F287").Value = "대한민국"
$ws.Range("G287").Value = "취업자 / 15세 이상 인구"
$ws.Range("H287").Value = "조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자"
$ws.Range("I287").Value = "평균 : 100만원, 중위값 : 1,000만원"
$ws.Range("J287").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("K287").Value = "`"19.1%`""
$ws.Range("L287").Value = "Red"
$ws.Range("M287").Value = "나는 사후 장기기증에 참여할 뜻이 없다"

$ws.Range("A288").Value = 45599.62022675926
$ws.Range("B288").Value = "chiyoon12@gmail.com"
$ws.Range("C288").Value = "경영학부"
$ws.Range("D288").Value = 20232938
$ws.Range("E288").Value = "김치윤"
$ws.Range("F288").Value = "대한민국"
$ws.Range("G288").Value = "취업자 / 경제활동인구"
$ws.Range("H288").Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Range("I288").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J288").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K288").Value = "`"10%`""
$ws.Range("L288").Value = "Black"
$ws.Range("N288").Value = "나는 사후 장기기증에 참여할 뜻이 있다"

$ws.Range("A289").Value = 45599.64794180555
$ws.Range("B289").Value = "blake4102@naver.com"
$ws.Range("C289").Value = "바이오메디컬"
$ws.Range("D289").Value = 20213609
$ws.Range("E289").Value = "김원래"
$ws.Range("F289").Value = "대한민국"
$ws.Range("G289").Value = "취업자 / 15세 이상 인구"
$ws.Range("H289").Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Range("I289").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J289").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K289").Value = "`"15%`""
$ws.Range("L289").Value = "Black"
$ws.Range("N289").Value = "나는 사후 장기기증에 참여할 뜻이 있다"

$ws.Range("A290").Value = 45599.64817425926
$ws.Range("B290").Value = "leeyelim0320@gmail.com"
$ws.Range("C290").Value = "바이오메디컬학과"
$ws.Range("D290").Value = 20193633
$ws.Range("E290").Value = "이예림"
$ws.Range("F290").Value = "미국"
$ws.Range("G290").Value = "경제활동인구 / 15세이상 인구"
$ws.Range("H290").Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Range("I290").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J290").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K290").Value = "`"19.1%`""
$ws.Range("L290").Value = "Black"
$ws.Range("N290").Value = "나는 사후 장기기증에 참여할 뜻이 있다"

$ws.Range("A291").Value = 45599.65252949074
$ws.Range("B291").Value = "4080jjh@gmail.com"
$ws.Range("C291").Value = "정치행정학과"
$ws.Range("D291").Value = 20182436
$ws.Range("E291").Value = "장재환"
$ws.Range("F291").Value = "미국"
$ws.Range("G291").Value = "경제활동인구 / 15세이상 인구"
$ws.Range("H291").Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Range("I291").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J291").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K291").Value = "`"19.1%`""
$ws.Range("L291").Value = "Red"
$ws.Range("M291").Value = "나는 사후 장기기증에 참여할 뜻이 없다"

$ws.Range("A292").Value = 45599.65902129629
$ws.Range("B292").Value = "minsung5342@naver.com"
$ws.Range("C292").Value = "사회복지학과"
$ws.Range("D292").Value = 20232311
$ws.Range("E292").Value = "김민성"
$ws.Range("F292").Value = "대한민국"
$ws.Range("G292").Value = "취업자 / 15세 이상 인구"
$ws.Range("H292").Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Range("I292").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J292").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K292").Value = "`"19.1%`""
$ws.Range("L292").Value = "Red"

$ws.Range("A293").Value = 45599.681940752314
$ws.Range("B293").Value = "youngold057@gmail.com"
$ws.Range("C293").Value = "사회복지학부"
$ws.Range("D293").Value = 20242336
$ws.Range("E293").Value = "윤태영"
$ws.Range("F293").Value = "대한민국"
$ws.Range("G293").Value = "취업자 / 15세 이상 인구"
$ws.Range("H293").Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Range("I293").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J293").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K293").Value = "`"19.1%`""
$ws.Range("L293").Value = "Red"

$ws.Range("A294").Value = 45599.68344836806
$ws.Range("B294").Value = "junhyeogjang@gmail.com"
$ws.Range("C294").Value = "반도체 디스플레이 스쿨"
$ws.Range("D294").Value = 20193341
$ws.Range("E294").Value = "장준혁"
$ws.Range("F294").Value = "미국"
$ws.Range("G294").Value = "실업자 / 경제활동인구"
$ws.Range("H294").Value = "조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자"
$ws.Range("I294").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("J294").Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Range("K294").Value = "`"19.1%`""
$ws.Range("L294").Value = "Black"

$ws.Range("A295").Value = 45599.68392090278
$ws.Range("B295").Value = "leedug87@gmail.com"
$ws.Range("C295").Value = "일본학과"
$ws.Range("D295").Value = 20231623
$ws.Range("E295").Value = "이두현"
$ws.Range("F295").Value = "대한민국"
$ws.Range("G295").Value = "경제활동인구 / 15세이상 인구"
$ws.Range("H295").Value = "조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자"
$ws.Range("I295").Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Range("J295").Value = "평균 : 100만원, 중위값 : 1,000만원"
$ws.Range("K295").Value = "`"19.1%`""
$ws.Range("L295").Value = "Black"

# ---------------------------------------------------------------------------
# 3. Row 283 is no longer the last row of the table, so it loses the two
#    trailing blank/placeholder cells (M283, N283) it used to carry.
# ---------------------------------------------------------------------------
$ws.Range("M283:N283").Clear()

# ---------------------------------------------------------------------------
# 3b. The format-copy step stamped every new row with both M and N (since
#     the even/odd template rows have both styled); drop whichever of the
#     two the real data doesn't occupy so the row ends exactly where the
#     source data does.
# ---------------------------------------------------------------------------
$ws.Range("M284").Clear()
$ws.Range("N285").Clear()
$ws.Range("M286:N286").Clear()
$ws.Range("N287").Clear()
$ws.Range("M288").Clear()
$ws.Range("M289").Clear()
$ws.Range("M290").Clear()
$ws.Range("N291").Clear()
$ws.Range("M292:N292").Clear()
$ws.Range("M293:N293").Clear()
$ws.Range("M294:N294").Clear()

# ---------------------------------------------------------------------------
# 4. Grow the "Form_Responses1" table to cover the 12 freshly added rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Form_Responses1")
$lo.Resize($ws.Range("A1:N295"))

# ---------------------------------------------------------------------------
# 5. Leave the cursor where the author's edit session left it.
# ---------------------------------------------------------------------------
$ws.Range("B300").Select()
